$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated simulation results (more games simulated, faster simulate-game logic,
# and refreshed optimization-derived transition probabilities).
$updates = @{
    "B2" = 0.2132867132867133
    "C2" = 0.5174825174825175
    "J2" = 0.01398601398601399
    "P2" = 0.1398601398601399
    "S2" = 0.1153846153846154
    "B3" = 0.00641025641025641
    "C3" = 0.03205128205128205
    "J3" = 0.02564102564102564
    "P3" = 0.7692307692307693
    "S3" = 0.1666666666666667
    "J4" = 0.05454545454545454
    "P4" = 0.7454545454545455
    "S4" = 0.2
    "B6" = 0.06880733944954129
    "D6" = 0.01376146788990826
    "F6" = 0.04587155963302753
    "J6" = 0.2201834862385321
    "O6" = 0.01834862385321101
    "Q6" = 0.1697247706422018
    "R6" = 0.08256880733944955
    "S6" = 0.3807339449541284
    "B7" = 0.1398963730569948
    "D7" = 0.0155440414507772
    "F7" = 0.04663212435233161
    "J7" = 0.1450777202072539
    "O7" = 0.03626943005181347
    "Q7" = 0.1295336787564767
    "R7" = 0.08808290155440414
    "S7" = 0.3989637305699482
    "B8" = 0.1035353535353535
    "D8" = 0.02777777777777778
    "F8" = 0.07323232323232323
    "J8" = 0.1237373737373737
    "O8" = 0.01767676767676768
    "Q8" = 0.154040404040404
    "R8" = 0.1161616161616162
    "S8" = 0.3838383838383838
    "B9" = 0.07801418439716312
    "D9" = 0.02127659574468085
    "E9" = 0.007092198581560284
    "F9" = 0.07801418439716312
    "J9" = 0.1063829787234043
    "O9" = 0.02127659574468085
    "Q9" = 0.1843971631205674
    "R9" = 0.07092198581560284
    "S9" = 0.4326241134751773
    "B10" = 0.1071428571428571
    "D10" = 0.02976190476190476
    "F10" = 0.09098639455782313
    "J10" = 0.1071428571428571
    "O10" = 0.03061224489795918
    "Q10" = 0.1743197278911565
    "R10" = 0.08163265306122448
    "S10" = 0.3784013605442177
    "G11" = 0.1446945337620579
    "J11" = 0.0932475884244373
    "K11" = 0.2057877813504823
    "L11" = 0.5466237942122186
    "S11" = 0.009646302250803859
    "G12" = 0.7100591715976331
    "J12" = 0.2189349112426036
    "K12" = 0.02366863905325444
    "L12" = 0.01183431952662722
    "S12" = 0.03550295857988166
    "G13" = 0.6739130434782609
    "J13" = 0.2173913043478261
    "S13" = 0.108695652173913
    "F15" = 0.0136986301369863
    "H15" = 0.1095890410958904
    "I15" = 0.0776255707762557
    "J15" = 0.2831050228310502
    "K15" = 0.0547945205479452
    "M15" = 0.0182648401826484
    "N15" = 0.0045662100456621
    "O15" = 0.0639269406392694
    "S15" = 0.3744292237442922
    "F16" = 0.01058201058201058
    "H16" = 0.1534391534391534
    "I16" = 0.06878306878306878
    "J16" = 0.3703703703703703
    "K16" = 0.1058201058201058
    "M16" = 0.02116402116402116
    "O16" = 0.0582010582010582
    "S16" = 0.2116402116402116
    "F17" = 0.01424501424501425
    "H17" = 0.1680911680911681
    "I17" = 0.06837606837606838
    "J17" = 0.4472934472934473
    "K17" = 0.07407407407407407
    "M17" = 0.01994301994301994
    "N17" = 0.002849002849002849
    "O17" = 0.03988603988603989
    "S17" = 0.1652421652421652
    "F18" = 0.0160427807486631
    "H18" = 0.213903743315508
    "I18" = 0.0481283422459893
    "J18" = 0.4064171122994653
    "K18" = 0.0855614973262032
    "M18" = 0.0106951871657754
    "O18" = 0.05882352941176471
    "S18" = 0.160427807486631
    "F19" = 0.01370967741935484
    "H19" = 0.2
    "I19" = 0.06290322580645161
    "J19" = 0.3774193548387097
    "K19" = 0.1346774193548387
    "M19" = 0.02580645161290323
    "N19" = 0.001612903225806452
    "O19" = 0.06854838709677419
    "S19" = 0.1153225806451613
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

Write-Host "Updated $($updates.Count) cells"
